# The user cut the contents of column C ("List Name 2 (PLACEMENT FEES)")
# and inserted them into column B, shifting the old (now-unused) column C
# out entirely. This mirrors Excel's "Cut" + "Insert Cut Cells" workflow,
# which is why the resulting selection lands on the whole of column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Cut() | Out-Null
$ws.Columns.Item(2).Insert() | Out-Null

# Leave the sheet selected the way Excel does after an "Insert Cut Cells"
# on a whole column.
$ws.Columns.Item(2).Select() | Out-Null
